$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential / "as of" date text (shared string used by A11)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.5032734799839308
$ws.Range("E2").Value = 0.01363107062127544

$ws.Range("D3").Value = 0.2428883953396222
$ws.Range("E3").Value = 0.01061249241964823

$ws.Range("D4").Value = 0.09507831182646373
$ws.Range("E4").Value = 0.01232965606748881

$ws.Range("D5").Value = 0.1027730709091098
$ws.Range("E5").Value = 0.0243948923194206

$ws.Range("D6").Value = 0.02954799239945446
$ws.Range("E6").Value = 0.02674753902754312

$ws.Range("D7").Value = 0.02643874954141899
$ws.Range("E7").Value = 0.0213511531252546

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.01447206235408882

# Restore sheet protection (sheet was protected before this edit)
$ws.Protect()
